$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): corrected Hydrogen demand, Non-metallic minerals entry
# no longer applies here - clear D3 but keep it as a blank cell (matching
# the blank style already used by its neighbour C3).
$ws.Range("B3").Value = 173631.5503952337
$ws.Range("D3").ClearContents()
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# Row 4 (Methanol): corrected Chemicals demand
$ws.Range("C4").Value = 0

# Row 5 (Ammonia): corrected Chemicals demand
$ws.Range("C5").Value = 2977.411704809432

# Row 7: "Other" is now "Biogas", with corrected Non-metallic minerals demand
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 367.201041958187

# Row 8 (new): a fresh "Other" row. Copy row 7's layout down first so the
# still-blank cells (B8, C8) come into existence the same way the rest of
# the sheet's blank cells did, then fill in the row label and the new
# Non-metallic minerals figure.
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 393.9265512588964
